# Applies the row-22/row-23 swap and appends new row 57, per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that must no longer be present after the row 22/23 swap ---
$ws.Cells.Item(22, 32).ClearContents()   # AF22 (Bestämningsmetod) moves off row 22
$ws.Cells.Item(23, 6).ClearContents()    # F23 (Artnamn) moves off row 23
$ws.Cells.Item(23, 11).ClearContents()   # K23 (Ålder-Stadium) moves off row 23
$ws.Cells.Item(23, 26).ClearContents()   # Z23 (Starttid) moves off row 23
$ws.Cells.Item(23, 28).ClearContents()   # AB23 (Sluttid) moves off row 23

# --- Row 22 ---
$ws.Cells.Item(22, 1).Value = 112390882
$ws.Cells.Item(22, 2).Value = 90800
$ws.Cells.Item(22, 3).Value = 'Ovaliderad'
$ws.Cells.Item(22, 4).Value = 'LC'
$ws.Cells.Item(22, 5).Value = 4364
$ws.Cells.Item(22, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(22, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(22, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(22, 11).Value = 'teleomorf'
$ws.Cells.Item(22, 16).Value = 'Kärmsjöbäckens naturreservat (Kärmsjöbäckens naturreservat), Ång'
$ws.Cells.Item(22, 17).Value = 583082
$ws.Cells.Item(22, 18).Value = 7086752
$ws.Cells.Item(22, 19).Value = 10
$ws.Cells.Item(22, 20).Value = 'Västernorrland'
$ws.Cells.Item(22, 21).Value = 'Sollefteå'
$ws.Cells.Item(22, 22).Value = 'Ångermanland'
$ws.Cells.Item(22, 23).Value = 'Junsele'
$ws.Cells.Item(22, 25).Value = "'2023-09-29"
$ws.Cells.Item(22, 26).Value = '11:58'
$ws.Cells.Item(22, 27).Value = "'2023-09-29"
$ws.Cells.Item(22, 28).Value = '11:58'
$ws.Cells.Item(22, 30).Value = $false
$ws.Cells.Item(22, 31).Value = $false
$ws.Cells.Item(22, 33).Value = $false
$ws.Cells.Item(22, 49).Value = 'Helena Thau'
$ws.Cells.Item(22, 50).Value = 'Helena Thau'

# --- Row 23 ---
$ws.Cells.Item(23, 1).Value = 112410309
$ws.Cells.Item(23, 2).Value = 90229
$ws.Cells.Item(23, 3).Value = 'Ovaliderad'
$ws.Cells.Item(23, 4).Value = 'NT'
$ws.Cells.Item(23, 5).Value = 757
$ws.Cells.Item(23, 7).Value = 'Hapalopilus aurantiacus'
$ws.Cells.Item(23, 8).Value = '(Rostk.) Bondartsev & Singer'
$ws.Cells.Item(23, 16).Value = 'Kärmsjöbäckens naturreservat, Ång'
$ws.Cells.Item(23, 17).Value = 583107
$ws.Cells.Item(23, 18).Value = 7086784
$ws.Cells.Item(23, 19).Value = 20
$ws.Cells.Item(23, 20).Value = 'Västernorrland'
$ws.Cells.Item(23, 21).Value = 'Sollefteå'
$ws.Cells.Item(23, 22).Value = 'Ångermanland'
$ws.Cells.Item(23, 23).Value = 'Junsele'
$ws.Cells.Item(23, 25).Value = "'2023-09-29"
$ws.Cells.Item(23, 27).Value = "'2023-09-29"
$ws.Cells.Item(23, 30).Value = $false
$ws.Cells.Item(23, 31).Value = $false
$ws.Cells.Item(23, 33).Value = $false
$ws.Cells.Item(23, 49).Value = 'Daniel Rutschman'
$ws.Cells.Item(23, 50).Value = 'Daniel Rutschman'

# --- Row 57 ---
$ws.Cells.Item(57, 1).Value = 112394707
$ws.Cells.Item(57, 2).Value = 89706
$ws.Cells.Item(57, 3).Value = 'Ovaliderad'
$ws.Cells.Item(57, 4).Value = 'VU'
$ws.Cells.Item(57, 5).Value = 1503
$ws.Cells.Item(57, 6).Value = 'Gräddporing'
$ws.Cells.Item(57, 7).Value = 'Sidera lenis'
$ws.Cells.Item(57, 8).Value = '(P.Karst.) Miettinen'
$ws.Cells.Item(57, 16).Value = 'Stor-Kärmsjön, Stor-Kärmsjön, Junsele s:n, Ång'
$ws.Cells.Item(57, 17).Value = 583090
$ws.Cells.Item(57, 18).Value = 7086458
$ws.Cells.Item(57, 19).Value = 25
$ws.Cells.Item(57, 20).Value = 'Västernorrland'
$ws.Cells.Item(57, 21).Value = 'Sollefteå'
$ws.Cells.Item(57, 22).Value = 'Ångermanland'
$ws.Cells.Item(57, 23).Value = 'Junsele'
$ws.Cells.Item(57, 25).Value = "'2023-09-29"
$ws.Cells.Item(57, 27).Value = "'2023-09-29"
$ws.Cells.Item(57, 30).Value = $false
$ws.Cells.Item(57, 31).Value = $false
$ws.Cells.Item(57, 33).Value = $false
$ws.Cells.Item(57, 49).Value = 'Daniel Rutschman'
$ws.Cells.Item(57, 50).Value = 'Daniel Rutschman'

# --- Reset number format/style on the date-text cells so they don't keep the
#     auto-applied date style (content must stay literal text, un-styled) ---
$ws.Cells.Item(22, 25).Style = "Normal"
$ws.Cells.Item(22, 27).Style = "Normal"
$ws.Cells.Item(23, 25).Style = "Normal"
$ws.Cells.Item(23, 27).Style = "Normal"
$ws.Cells.Item(57, 25).Style = "Normal"
$ws.Cells.Item(57, 27).Style = "Normal"
